# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'63.395.52"
$ws.Range('E2').Value = '  -1.37%  '

# Row 3
$ws.Range('D3').Value = "'2.598.31"
$ws.Range('E3').Value = '  -1.23%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').Value = "'588.37"
$ws.Range('E5').Value = '  -2.71%  '

# Row 6
$ws.Range('D6').Value = "'149.11"
$ws.Range('E6').Value = '  -1.51%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').Value = "'0.583"
$ws.Range('E8').Value = '  -1.45%  '

# Row 9
$ws.Range('E9').Value = '  -0.97%  '

# Row 10
$ws.Range('D10').Value = "'5.73"
$ws.Range('E10').Value = '  -0.59%  '

# Row 11
$ws.Range('D11').Value = "'0.385"
$ws.Range('E11').Value = '  -0.31%  '

# Row 12
$ws.Range('D12').Value = "'0.151"
$ws.Range('E12').Value = '  -0.03%  '

# Row 13
$ws.Range('D13').Value = "'27.53"
$ws.Range('E13').Value = '  -0.55%  '

# Row 14
$ws.Range('D14').Value = "'3.068.27"
$ws.Range('E14').Value = '  -1.28%  '

# Row 15
$ws.Range('D15').Value = "'63.287.44"
$ws.Range('E15').Value = '  -1.34%  '

# Row 16
$ws.Range('D16').Value = "'0.0000156"
$ws.Range('E16').Value = '  +3.62%  '

# Row 17
$ws.Range('D17').Value = "'2.590.07"
$ws.Range('E17').Value = '  -1.70%  '

# Row 18
$ws.Range('D18').Value = "'12.07"
$ws.Range('E18').Value = '  -0.88%  '

# Row 19
$ws.Range('E19').Value = '  +0.50%  '

# Row 20
$ws.Range('D20').Value = "'343.90"
$ws.Range('E20').Value = '  -1.75%  '

# Row 21
$ws.Range('D21').Value = "'6.83"
$ws.Range('E21').Value = '  -2.34%  '

# Row 22
$ws.Range('E22').Value = '  +0.06%  '

# Row 23
$ws.Range('D23').Value = "'66.57"
$ws.Range('E23').Value = '  -0.31%  '

# Row 24
$ws.Range('D24').Value = "'1.71"
$ws.Range('E24').Value = '  -2.20%  '

# Row 25
$ws.Range('D25').Value = "'9.15"
$ws.Range('E25').Value = '  -1.52%  '

# Row 26
$ws.Range('D26').Value = "'1.64"
$ws.Range('E26').Value = '  -3.79%  '

# Row 27
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = "'8.25"
$ws.Range('E27').Value = '  +1.51%  '

# Row 28
$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D28').Value = "'554.46"
$ws.Range('E28').Value = '  +1.85%  '

# Row 29
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.01%  '

# Row 30
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = "'0.161"
$ws.Range('E30').Value = '  -3.68%  '

# Row 31
$ws.Range('E31').Value = '  -2.36%  '

# Row 32
$ws.Range('D32').Value = "'0.0₃0854"
$ws.Range('E32').Value = '  -0.83%  '

# Row 33
$ws.Range('E33').Value = '  -0.69%  '

# Row 34
$ws.Range('D34').Value = "'5.26"
$ws.Range('E34').Value = '  -0.44%  '

# Row 35
$ws.Range('D35').Value = "'6.06"
$ws.Range('E35').Value = '  -1.60%  '

# Row 36
$ws.Range('D36').Value = "'165.45"
$ws.Range('E36').Value = '  -1.20%  '

# Row 37
$ws.Range('D37').Value = "'0.412"
$ws.Range('E37').Value = '  +0.05%  '

# Row 38
$ws.Range('E38').Value = '  -0.06%  '

# Row 39
$ws.Range('E39').Value = '  -0.87%  '

# Row 40
$ws.Range('D40').Value = "'1.91"
$ws.Range('E40').Value = '  -5.83%  '

# Row 41
$ws.Range('E41').Value = '  -0.05%  '

# Row 42
$ws.Range('D42').Value = "'165.23"
$ws.Range('E42').Value = '  -3.81%  '

# Row 43
$ws.Range('D43').Value = "'3.99"
$ws.Range('E43').Value = '  +1.19%  '

# Row 44
$ws.Range('D44').Value = "'22.76"
$ws.Range('E44').Value = '  +5.26%  '

# Row 45
$ws.Range('E45').Value = '  -1.74%  '

# Row 46
$ws.Range('E46').Value = '  +4.17%  '

# Row 47
$ws.Range('D47').Value = "'0.630"
$ws.Range('E47').Value = '  +0.18%  '

# Row 48
$ws.Range('E48').Value = '  +0.13%  '

# Row 49
$ws.Range('D49').Value = "'0.0957"
$ws.Range('E49').Value = '  -1.02%  '

# Row 50
$ws.Range('D50').Value = "'19.00"
$ws.Range('E50').Value = '  -1.95%  '
